$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before the existing row 277, shifting old rows 277-303 down to 279-305.
$ws.Rows("277:278").Insert()

# Fill the 2 newly inserted rows (277-278) with the new week of price data.
$ws.Cells.Item(277,1).Value = 5
$ws.Cells.Item(277,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(277,3).Value = 'Maule'
$ws.Cells.Item(277,4).Value = 44449
$ws.Cells.Item(277,5).Value = 7
$ws.Cells.Item(277,6).Value = 100112002
$ws.Cells.Item(277,7).Value = 'Pimiento'
$ws.Cells.Item(277,8).Value = 'Zafiro rojo'
$ws.Cells.Item(277,9).Value = 'Primera'
$ws.Cells.Item(277,10).Value = 200
$ws.Cells.Item(277,11).Value = 43000
$ws.Cells.Item(277,12).Value = 43000
$ws.Cells.Item(277,13).Value = 43000
$ws.Cells.Item(277,14).Value = '$/caja 15 kilos'
$ws.Cells.Item(277,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(277,16).Value = 2867
$ws.Cells.Item(277,17).Value = 15
$ws.Cells.Item(277,18).Value = 'Hortaliza'

$ws.Cells.Item(278,1).Value = 5
$ws.Cells.Item(278,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(278,3).Value = 'Maule'
$ws.Cells.Item(278,4).Value = 44449
$ws.Cells.Item(278,5).Value = 7
$ws.Cells.Item(278,6).Value = 100112002
$ws.Cells.Item(278,7).Value = 'Pimiento'
$ws.Cells.Item(278,8).Value = 'Zafiro verde'
$ws.Cells.Item(278,9).Value = 'Primera'
$ws.Cells.Item(278,10).Value = 200
$ws.Cells.Item(278,11).Value = 40000
$ws.Cells.Item(278,12).Value = 40000
$ws.Cells.Item(278,13).Value = 40000
$ws.Cells.Item(278,14).Value = '$/caja 15 kilos'
$ws.Cells.Item(278,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(278,16).Value = 2667
$ws.Cells.Item(278,17).Value = 15
$ws.Cells.Item(278,18).Value = 'Hortaliza'

# Append 2 new rows (304-305) at the end with the rolled-over data.
$ws.Cells.Item(304,1).Value = 5
$ws.Cells.Item(304,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(304,3).Value = 'Maule'
$ws.Cells.Item(304,4).Value = 44400
$ws.Cells.Item(304,5).Value = 7
$ws.Cells.Item(304,6).Value = 100112002
$ws.Cells.Item(304,7).Value = 'Pimiento'
$ws.Cells.Item(304,8).Value = 'Zafiro rojo'
$ws.Cells.Item(304,9).Value = 'Primera'
$ws.Cells.Item(304,10).Value = 300
$ws.Cells.Item(304,11).Value = 17000
$ws.Cells.Item(304,12).Value = 17000
$ws.Cells.Item(304,13).Value = 17000
$ws.Cells.Item(304,14).Value = '$/caja 15 kilos'
$ws.Cells.Item(304,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(304,16).Value = 1133
$ws.Cells.Item(304,17).Value = 15
$ws.Cells.Item(304,18).Value = 'Hortaliza'
$ws.Cells.Item(304,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(305,1).Value = 5
$ws.Cells.Item(305,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(305,3).Value = 'Maule'
$ws.Cells.Item(305,4).Value = 44400
$ws.Cells.Item(305,5).Value = 7
$ws.Cells.Item(305,6).Value = 100112002
$ws.Cells.Item(305,7).Value = 'Pimiento'
$ws.Cells.Item(305,8).Value = 'Zafiro verde'
$ws.Cells.Item(305,9).Value = 'Primera'
$ws.Cells.Item(305,10).Value = 300
$ws.Cells.Item(305,11).Value = 15000
$ws.Cells.Item(305,12).Value = 15000
$ws.Cells.Item(305,13).Value = 15000
$ws.Cells.Item(305,14).Value = '$/caja 15 kilos'
$ws.Cells.Item(305,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(305,16).Value = 1000
$ws.Cells.Item(305,17).Value = 15
$ws.Cells.Item(305,18).Value = 'Hortaliza'
$ws.Cells.Item(305,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"